$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab to "CubeA"
$ws.Name = "CubeA"

# Add new row 16 with data, mirroring the style pattern of previous rows
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(16, 1).Value = 14

$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item(16, 3).Value = 0.9671672344821337
$ws.Cells.Item(16, 4).Value = 1.044982739407209
$ws.Cells.Item(16, 5).Value = 0.985526333062672
$ws.Cells.Item(16, 6).Value = 1.006797359985314
$ws.Cells.Item(16, 7).Value = 0.9671672344821337
$ws.Cells.Item(16, 8).Value = 1.044982739407209
$ws.Cells.Item(16, 9).Value = 0.9873632744601039
$ws.Cells.Item(16, 10).Value = 1.007414843407916
$ws.Cells.Item(16, 11).Value = 0.9877485424670935
$ws.Cells.Item(16, 12).Value = 1.023711712484567
$ws.Cells.Item(16, 13).Value = 0.9671672344821337
$ws.Cells.Item(16, 14).Value = 1.015254536234941
$ws.Cells.Item(16, 15).Value = 1.001118416734332
$ws.Cells.Item(16, 16).Value = 1.001339004969626
